$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B17:F18").Value = 0.2

$ws.Range("E28").Value = "Abnahme/ Übergabe; Projektstrukturierung; Meilensteinplanung; Qualitätssicherung"
$ws.Range("D28").Value = "Einleitung; Analyse der Infrage kommenden Produkte; Wirtschaftlichkeits-analyse "
$ws.Range("F28").Value = "Umsetzung; Serverinstallation und -Konfiguration;  Clientinstallation und -Konfiguration"
$ws.Range("B28").Value = "Projekt- /Prozessschnittstellen; IST-Analyse; SOLL-Analyse"

$ws.Range("B29").Value = "Definitionsphase"
$ws.Range("C29").Value = "Planungsphase"
$ws.Range("E29").Value = "Abschlussphase"
$ws.Range("F29").Value = "Durchführungsphase"
$ws.Range("D29").Value = "Format; zusammenfügen der Phasen; Animationen; Überarbeitung"

$ws.Range("C28").Value = "Projektumfeld; Projektziel; Entscheidung für ein Produkt; Testfallkatalog; Risikoanalyse; Zeitplanung; Fazit"

$ws.Rows.Item(28).RowHeight = 90

$ws.Range("D19").Select()

$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 2
